# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (used only by the notes master)
#   ppt/theme/theme2.xml  -> "Integral"     (used by the slide master / all slides)
#
# The authored commit swaps the contents of these two theme parts: theme1.xml
# becomes the "Integral" theme and theme2.xml becomes the plain "Office Theme"
# (font scheme and format scheme are identical between the two parts already,
# only the colour scheme - and its name - actually differ).
#
# The only theme surface the PowerPoint object model exposes for writing is
# the live/active theme that backs the slide master (reached here via any
# slide's ThemeColorScheme); the notes-master's theme1.xml is not reachable
# for writes through the object model, so we repoint the live theme (which
# lands in theme2.xml) to the "Office Theme" colours.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Theme colour scheme slot order is fixed: 1=dk1 2=lt1 3=dk2 4=lt2
# 5=accent1 6=accent2 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
# Values are plain VBA RGB() integers (r + g*256 + b*65536) for the "Office
# Theme" palette that used to live in theme1.xml.
$cs.Item(1).RGB  = 0         # dk1      000000
$cs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      44546A
$cs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  FFC000
$cs.Item(9).RGB  = 12874308  # accent5  4472C4
$cs.Item(10).RGB = 4697456   # accent6  70AD47
$cs.Item(11).RGB = 12673797  # hlink    0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72
